# edit.ps1
# Applies the "Doing Updates for Financials" revision to JBHT_QTR_FIN.xlsx:
#  - Inserts two new quarterly columns (D,E) ahead of the existing data,
#    shifting the prior D:K data block to F:M.
#  - Populates the two new columns with the newest two quarters of data
#    for the Income Statement, Balance Sheet and Cash Flow Statement blocks.
#  - Applies a handful of follow-up corrections to the quarter that lands in
#    column H (previously column D->F... originally column D) for a few line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert two blank columns at D:E. This pushes the existing quarterly
#    columns D:K to F:M, exactly matching the diff column remapping.
# ---------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# ---------------------------------------------------------------------
# 2) The newly inserted columns inherit the "General" format; copy the
#    number format from column F (which holds the old column D data)
#    into the new D:E columns so dates/numbers render the same way.
# ---------------------------------------------------------------------
$dataRows = @(7,8,9,10,12,13,14,15,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,38,41,42,43,44,45,46,47,48,49,50,51,52,53,54,57,58,59,60,61,62,63,64,65,66,68,69,70,71,72,73,74,75,76,77,80,81,83,84,85,86,87,88,89,91,92,93,94,96,97,98,99,100,101,102)
foreach ($r in $dataRows) {
    $fmt = $ws.Cells.Item($r, 6).NumberFormat
    $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 5)).NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 3) Populate the two new columns (D = newest quarter, E = prior quarter)
#    for every row that carries data.
# ---------------------------------------------------------------------
$ws.Cells.Item(7,4).Value2 = 43465
$ws.Cells.Item(7,5).Value2 = 43373
$ws.Cells.Item(8,4).Value2 = 2317800
$ws.Cells.Item(8,5).Value2 = 2209800
$ws.Cells.Item(9,4).Value2 = 1981200
$ws.Cells.Item(9,5).Value2 = 1817900
$ws.Cells.Item(10,4).Value2 = 336600
$ws.Cells.Item(10,5).Value2 = 391900
$ws.Cells.Item(12,4).Value2 = "NA"
$ws.Cells.Item(12,5).Value2 = "NA"
$ws.Cells.Item(13,4).Value2 = 0
$ws.Cells.Item(13,5).Value2 = 0
$ws.Cells.Item(14,4).Value2 = 0
$ws.Cells.Item(14,5).Value2 = 0
$ws.Cells.Item(15,4).Value2 = 114100
$ws.Cells.Item(15,5).Value2 = 108800
$ws.Cells.Item(17,4).Value2 = 2195100
$ws.Cells.Item(17,5).Value2 = 2035100
$ws.Cells.Item(18,4).Value2 = 122700
$ws.Cells.Item(18,5).Value2 = 174700
$ws.Cells.Item(20,4).Value2 = -11200
$ws.Cells.Item(20,5).Value2 = -10000
$ws.Cells.Item(21,4).Value2 = 225600
$ws.Cells.Item(21,5).Value2 = 273500
$ws.Cells.Item(22,4).Value2 = 0
$ws.Cells.Item(22,5).Value2 = 0
$ws.Cells.Item(23,4).Value2 = 111500
$ws.Cells.Item(23,5).Value2 = 164700
$ws.Cells.Item(24,4).Value2 = 6000
$ws.Cells.Item(24,5).Value2 = 33600
$ws.Cells.Item(25,4).Value2 = 0
$ws.Cells.Item(25,5).Value2 = 0
$ws.Cells.Item(26,4).Value2 = 105500
$ws.Cells.Item(26,5).Value2 = 131100
$ws.Cells.Item(27,4).Value2 = 105500
$ws.Cells.Item(27,5).Value2 = 131100
$ws.Cells.Item(28,4).Value2 = 0
$ws.Cells.Item(28,5).Value2 = 0
$ws.Cells.Item(29,4).Value2 = -16800
$ws.Cells.Item(29,5).Value2 = 0
$ws.Cells.Item(30,4).Value2 = 0
$ws.Cells.Item(30,5).Value2 = 0
$ws.Cells.Item(31,4).Value2 = 0
$ws.Cells.Item(31,5).Value2 = 0
$ws.Cells.Item(32,4).Value2 = 11200
$ws.Cells.Item(32,5).Value2 = 10000
$ws.Cells.Item(33,4).Value2 = 88700
$ws.Cells.Item(33,5).Value2 = 131100
$ws.Cells.Item(34,4).Value2 = 0
$ws.Cells.Item(34,5).Value2 = 0
$ws.Cells.Item(35,4).Value2 = 88700
$ws.Cells.Item(35,5).Value2 = 131100
$ws.Cells.Item(38,4).Value2 = 43465
$ws.Cells.Item(38,5).Value2 = 43373
$ws.Cells.Item(41,4).Value2 = 7600
$ws.Cells.Item(41,5).Value2 = 7600
$ws.Cells.Item(42,4).Value2 = 0
$ws.Cells.Item(42,5).Value2 = 0
$ws.Cells.Item(43,4).Value2 = 1326200
$ws.Cells.Item(43,5).Value2 = 1095600
$ws.Cells.Item(44,4).Value2 = 22000
$ws.Cells.Item(44,5).Value2 = "NA"
$ws.Cells.Item(45,4).Value2 = 147200
$ws.Cells.Item(45,5).Value2 = 239300
$ws.Cells.Item(46,4).Value2 = 1503000
$ws.Cells.Item(46,5).Value2 = 1342500
$ws.Cells.Item(47,4).Value2 = "NA"
$ws.Cells.Item(47,5).Value2 = 17500
$ws.Cells.Item(48,4).Value2 = 3445100
$ws.Cells.Item(48,5).Value2 = 3241500
$ws.Cells.Item(49,4).Value2 = 105200
$ws.Cells.Item(49,5).Value2 = "NA"
$ws.Cells.Item(50,4).Value2 = 0
$ws.Cells.Item(50,5).Value2 = 0
$ws.Cells.Item(51,4).Value2 = 0
$ws.Cells.Item(51,5).Value2 = 0
$ws.Cells.Item(52,4).Value2 = 38400
$ws.Cells.Item(52,5).Value2 = 120200
$ws.Cells.Item(53,4).Value2 = 0
$ws.Cells.Item(53,5).Value2 = 0
$ws.Cells.Item(54,4).Value2 = 5091600
$ws.Cells.Item(54,5).Value2 = 4721700
$ws.Cells.Item(57,4).Value2 = 709700
$ws.Cells.Item(57,5).Value2 = 569500
$ws.Cells.Item(58,4).Value2 = 250700
$ws.Cells.Item(58,5).Value2 = 248700
$ws.Cells.Item(59,4).Value2 = 391900
$ws.Cells.Item(59,5).Value2 = 394500
$ws.Cells.Item(60,4).Value2 = 1352300
$ws.Cells.Item(60,5).Value2 = 1212700
$ws.Cells.Item(61,4).Value2 = 898400
$ws.Cells.Item(61,5).Value2 = 820900
$ws.Cells.Item(62,4).Value2 = 739500
$ws.Cells.Item(62,5).Value2 = 607600
$ws.Cells.Item(63,4).Value2 = 0
$ws.Cells.Item(63,5).Value2 = 0
$ws.Cells.Item(64,4).Value2 = 0
$ws.Cells.Item(64,5).Value2 = 0
$ws.Cells.Item(65,4).Value2 = 0
$ws.Cells.Item(65,5).Value2 = 0
$ws.Cells.Item(66,4).Value2 = 2990300
$ws.Cells.Item(66,5).Value2 = 2641100
$ws.Cells.Item(68,4).Value2 = 0
$ws.Cells.Item(68,5).Value2 = 0
$ws.Cells.Item(69,4).Value2 = 0
$ws.Cells.Item(69,5).Value2 = 0
$ws.Cells.Item(70,4).Value2 = 0
$ws.Cells.Item(70,5).Value2 = 0
$ws.Cells.Item(71,4).Value2 = 0
$ws.Cells.Item(71,5).Value2 = 0
$ws.Cells.Item(72,4).Value2 = 4188400
$ws.Cells.Item(72,5).Value2 = "NA"
$ws.Cells.Item(73,4).Value2 = 0
$ws.Cells.Item(73,5).Value2 = 0
$ws.Cells.Item(74,4).Value2 = 0
$ws.Cells.Item(74,5).Value2 = 0
$ws.Cells.Item(75,4).Value2 = 0
$ws.Cells.Item(75,5).Value2 = 0
$ws.Cells.Item(76,4).Value2 = 2101400
$ws.Cells.Item(76,5).Value2 = 2080600
$ws.Cells.Item(77,4).Value2 = 0
$ws.Cells.Item(77,5).Value2 = 0
$ws.Cells.Item(80,4).Value2 = 43465
$ws.Cells.Item(80,5).Value2 = 43373
$ws.Cells.Item(81,4).Value2 = 88700
$ws.Cells.Item(81,5).Value2 = 131100
$ws.Cells.Item(83,4).Value2 = 114100
$ws.Cells.Item(83,5).Value2 = 108800
$ws.Cells.Item(84,4).Value2 = 0
$ws.Cells.Item(84,5).Value2 = 0
$ws.Cells.Item(85,4).Value2 = 0
$ws.Cells.Item(85,5).Value2 = 0
$ws.Cells.Item(86,4).Value2 = 0
$ws.Cells.Item(86,5).Value2 = 0
$ws.Cells.Item(87,4).Value2 = 0
$ws.Cells.Item(87,5).Value2 = 0
$ws.Cells.Item(88,4).Value2 = 0
$ws.Cells.Item(88,5).Value2 = 0
$ws.Cells.Item(89,4).Value2 = 310200
$ws.Cells.Item(89,5).Value2 = 248600
$ws.Cells.Item(91,4).Value2 = -319600
$ws.Cells.Item(91,5).Value2 = -264400
$ws.Cells.Item(92,4).Value2 = 0
$ws.Cells.Item(92,5).Value2 = 0
$ws.Cells.Item(93,4).Value2 = 0
$ws.Cells.Item(93,5).Value2 = 0
$ws.Cells.Item(94,4).Value2 = -300500
$ws.Cells.Item(94,5).Value2 = -230800
$ws.Cells.Item(96,4).Value2 = -26200
$ws.Cells.Item(96,5).Value2 = -26200
$ws.Cells.Item(97,4).Value2 = 0
$ws.Cells.Item(97,5).Value2 = 0
$ws.Cells.Item(98,4).Value2 = 0
$ws.Cells.Item(98,5).Value2 = 0
$ws.Cells.Item(99,4).Value2 = 0
$ws.Cells.Item(99,5).Value2 = 0
$ws.Cells.Item(100,4).Value2 = -9700
$ws.Cells.Item(100,5).Value2 = -25300
$ws.Cells.Item(101,4).Value2 = 0
$ws.Cells.Item(101,5).Value2 = 0
$ws.Cells.Item(102,4).Value2 = 0
$ws.Cells.Item(102,5).Value2 = -7600

# ---------------------------------------------------------------------
# 4) A few line items also received a data correction in the quarter that
#    now sits in column H (previously column D, then shifted to F... the
#    original "Jun-17" column), beyond the pure column shift above.
# ---------------------------------------------------------------------
$ws.Cells.Item(24,8).Value2 = -554700
$ws.Cells.Item(26,8).Value2 = 694500
$ws.Cells.Item(27,8).Value2 = 694500
$ws.Cells.Item(29,8).Value2 = -309200

Write-Host "Done"
